$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column L ("Estado") -- several transactions previously marked as
# FAILED / REJECTED are now APPROVED (physical file generation succeeded).
$ws.Range("L6").Value = "APPROVED"
$ws.Range("L7").Value = "APPROVED"
$ws.Range("L8").Value = "APPROVED"
$ws.Range("L9").Value = "APPROVED"
$ws.Range("L12").Value = "APPROVED"

# Column M ("Monto bruto") -- gross amount corrected/populated with the
# real transaction amounts now that the transactions are approved.
$ws.Range("M8").Value = 40008.54
$ws.Range("M10").Value = 16846.35
$ws.Range("M11").Value = 1111235.55
$ws.Range("M13").Value = 115.25
$ws.Range("M14").Value = 130.48
